$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 649.9
$ws.Cells.Item(4, 9).Value = 547
$ws.Cells.Item(4, 11).Value = 547
$ws.Cells.Item(4, 13).Value = -433

$ws.Cells.Item(21, 8).Value = 2833
$ws.Cells.Item(21, 10).Value = 3000
$ws.Cells.Item(21, 12).Value = 3000
$ws.Cells.Item(21, 14).Value = -3936

$ws.Cells.Item(23, 8).Value = 2833
$ws.Cells.Item(23, 10).Value = 3000
$ws.Cells.Item(23, 12).Value = 3000
$ws.Cells.Item(23, 14).Value = -3468

$ws.Cells.Item(37, 8).Value = 1400
$ws.Cells.Item(37, 9).Value = 1400
$ws.Cells.Item(37, 11).Value = 4200
$ws.Cells.Item(37, 13).Value = -4074

$ws.Cells.Item(40, 8).Value = 13263.048
$ws.Cells.Item(40, 9).Value = 7353.8335
$ws.Cells.Item(40, 11).Value = 7353.8335
$ws.Cells.Item(40, 13).Value = -7178.8335

$ws.Cells.Item(51, 8).Value = 3760
$ws.Cells.Item(51, 10).Value = 3950
$ws.Cells.Item(51, 12).Value = 3950
$ws.Cells.Item(51, 14).Value = -4918

$ws.Cells.Item(53, 8).Value = 1599.4166
$ws.Cells.Item(53, 9).Value = 1466.6666
$ws.Cells.Item(53, 10).Value = 1732.1666
$ws.Cells.Item(53, 11).Value = 1466.6666
$ws.Cells.Item(53, 12).Value = 1732.1666
$ws.Cells.Item(53, 13).Value = -829.6666
$ws.Cells.Item(53, 14).Value = -3006.1666

$ws.Cells.Item(132, 8).Value = 1739.0834
$ws.Cells.Item(132, 9).Value = 1771.2174
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 5313.6522
$ws.Cells.Item(132, 12).Value = 3000
$ws.Cells.Item(132, 13).Value = -2783.6522
$ws.Cells.Item(132, 14).Value = -8060

$ws.Cells.Item(135, 8).Value = 1229.2307
$ws.Cells.Item(135, 10).Value = 941.8
$ws.Cells.Item(135, 12).Value = 8476.199999999999
$ws.Cells.Item(135, 14).Value = -13546.2

$ws.Cells.Item(137, 8).Value = 299668.3
$ws.Cells.Item(137, 9).Value = 3509.639
$ws.Cells.Item(137, 10).Value = 1119800.1
$ws.Cells.Item(137, 11).Value = 10528.917
$ws.Cells.Item(137, 12).Value = 3359400.3
$ws.Cells.Item(137, 13).Value = -7978.917000000001
$ws.Cells.Item(137, 14).Value = -3364500.3

$ws.Cells.Item(138, 8).Value = 2644.102
$ws.Cells.Item(138, 9).Value = 2065.625
$ws.Cells.Item(138, 10).Value = 2924.5757
$ws.Cells.Item(138, 11).Value = 6196.875
$ws.Cells.Item(138, 12).Value = 8773.7271
$ws.Cells.Item(138, 13).Value = -1056.875
$ws.Cells.Item(138, 14).Value = -19053.7271

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 1233.4
$ws.Cells.Item(35, 10).Value = 959.3333
$ws.Cells.Item(35, 12).Value = 959.3333
$ws.Cells.Item(35, 14).Value = -1771.3333

$ws.Cells.Item(51, 8).Value = 32000
$ws.Cells.Item(51, 10).Value = 32000
$ws.Cells.Item(51, 12).Value = 32000
$ws.Cells.Item(51, 14).Value = -33512

$ws.Cells.Item(61, 8).Value = 55118.79
$ws.Cells.Item(61, 9).Value = 2508.5386
$ws.Cells.Item(61, 10).Value = 169107.67
$ws.Cells.Item(61, 11).Value = 2508.5386
$ws.Cells.Item(61, 12).Value = 169107.67
$ws.Cells.Item(61, 13).Value = -2296.5386
$ws.Cells.Item(61, 14).Value = -169531.67

$ws.Cells.Item(63, 8).Value = 5064.143
$ws.Cells.Item(63, 9).Value = 1889.8
$ws.Cells.Item(63, 11).Value = 1889.8
$ws.Cells.Item(63, 13).Value = -1203.8

$ws.Cells.Item(66, 8).Value = 5064.143
$ws.Cells.Item(66, 9).Value = 1889.8
$ws.Cells.Item(66, 11).Value = 9449
$ws.Cells.Item(66, 13).Value = -6017

$ws.Cells.Item(136, 8).Value = 55118.79
$ws.Cells.Item(136, 9).Value = 2508.5386
$ws.Cells.Item(136, 10).Value = 169107.67
$ws.Cells.Item(136, 11).Value = 7525.6158
$ws.Cells.Item(136, 12).Value = 507323.01
$ws.Cells.Item(136, 13).Value = -4975.6158
$ws.Cells.Item(136, 14).Value = -512423.01

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 198934.77
$ws.Cells.Item(20, 9).Value = 272413.8
$ws.Cells.Item(20, 11).Value = 272413.8
$ws.Cells.Item(20, 13).Value = -272166.8

$ws.Cells.Item(107, 8).Value = 2988.7368
$ws.Cells.Item(107, 9).Value = 2259.182
$ws.Cells.Item(107, 11).Value = 2259.182
$ws.Cells.Item(107, 13).Value = -339.1819999999998

$ws.Cells.Item(134, 8).Value = 1926.3448
$ws.Cells.Item(134, 9).Value = 960.4737
$ws.Cells.Item(134, 11).Value = 2881.4211
$ws.Cells.Item(134, 13).Value = -346.4211

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1408
$ws.Cells.Item(16, 9).Value = 1071.2
$ws.Cells.Item(16, 11).Value = 1071.2
$ws.Cells.Item(16, 13).Value = -784.2

$ws.Cells.Item(31, 8).Value = 3108.8572
$ws.Cells.Item(31, 9).Value = 2685.3333
$ws.Cells.Item(31, 10).Value = 3871.2
$ws.Cells.Item(31, 11).Value = 2685.3333
$ws.Cells.Item(31, 12).Value = 3871.2
$ws.Cells.Item(31, 13).Value = -2390.3333
$ws.Cells.Item(31, 14).Value = -4461.2

$ws.Cells.Item(34, 8).Value = 3108.8572
$ws.Cells.Item(34, 9).Value = 2685.3333
$ws.Cells.Item(34, 10).Value = 3871.2
$ws.Cells.Item(34, 11).Value = 2685.3333
$ws.Cells.Item(34, 12).Value = 3871.2
$ws.Cells.Item(34, 13).Value = -2483.3333
$ws.Cells.Item(34, 14).Value = -4275.2

$ws.Cells.Item(58, 8).Value = 2116.36
$ws.Cells.Item(58, 9).Value = 1888.8
$ws.Cells.Item(58, 10).Value = 2268.0667
$ws.Cells.Item(58, 11).Value = 1888.8
$ws.Cells.Item(58, 12).Value = 2268.0667
$ws.Cells.Item(58, 13).Value = -1685.8
$ws.Cells.Item(58, 14).Value = -2674.0667

$ws.Cells.Item(107, 8).Value = 187.42857
$ws.Cells.Item(107, 9).Value = 182.4
$ws.Cells.Item(107, 11).Value = 182.4
$ws.Cells.Item(107, 13).Value = 1737.6

$ws.Cells.Item(113, 8).Value = 1408
$ws.Cells.Item(113, 9).Value = 1071.2
$ws.Cells.Item(113, 11).Value = 1071.2
$ws.Cells.Item(113, 13).Value = 1098.8

$ws.Cells.Item(134, 8).Value = 38156.586
$ws.Cells.Item(134, 9).Value = 4277.15
$ws.Cells.Item(134, 10).Value = 113444.22
$ws.Cells.Item(134, 11).Value = 12831.45
$ws.Cells.Item(134, 12).Value = 340332.66
$ws.Cells.Item(134, 13).Value = -10296.45
$ws.Cells.Item(134, 14).Value = -345402.66

$ws.Cells.Item(136, 8).Value = 2116.36
$ws.Cells.Item(136, 9).Value = 1888.8
$ws.Cells.Item(136, 10).Value = 2268.0667
$ws.Cells.Item(136, 11).Value = 5666.4
$ws.Cells.Item(136, 12).Value = 6804.2001
$ws.Cells.Item(136, 13).Value = -3116.4
$ws.Cells.Item(136, 14).Value = -11904.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(104, 8).Value = 4541.0527
$ws.Cells.Item(104, 10).Value = 5000
$ws.Cells.Item(104, 12).Value = 15000
$ws.Cells.Item(104, 14).Value = -20242

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 2015.1428
$ws.Cells.Item(2, 9).Value = 3833.6667
$ws.Cells.Item(2, 10).Value = 651.25
$ws.Cells.Item(2, 11).Value = 3833.6667
$ws.Cells.Item(2, 12).Value = 651.25
$ws.Cells.Item(2, 13).Value = -3720.6667
$ws.Cells.Item(2, 14).Value = -877.25

$ws.Cells.Item(107, 8).Value = 276.10526
$ws.Cells.Item(107, 9).Value = 101.111115
$ws.Cells.Item(107, 11).Value = 101.111115
$ws.Cells.Item(107, 13).Value = 1818.888885

$ws.Cells.Item(132, 8).Value = 5374.5713
$ws.Cells.Item(132, 9).Value = 3547.7917
$ws.Cells.Item(132, 11).Value = 10643.3751
$ws.Cells.Item(132, 13).Value = -8113.375100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).Value = ""

$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = ""
$ws.Cells.Item(61, 14).Value = ""

$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = ""
$ws.Cells.Item(113, 14).Value = ""

$ws.Cells.Item(122, 8).Value = 11788141
$ws.Cells.Item(122, 10).Value = 33337966
$ws.Cells.Item(122, 12).Value = 100013898
$ws.Cells.Item(122, 14).Value = -100018798

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3107504
$ws.Cells.Item(132, 9).Value = 606.6667
$ws.Cells.Item(132, 10).Value = 8699919
$ws.Cells.Item(132, 11).Value = 1820.0001
$ws.Cells.Item(132, 12).Value = 26099757
$ws.Cells.Item(132, 13).Value = 709.9999
$ws.Cells.Item(132, 14).Value = -26104817

$ws.Cells.Item(136, 8).Value = 1706.4166
$ws.Cells.Item(136, 10).Value = 1848.1538
$ws.Cells.Item(136, 12).Value = 5544.4614
$ws.Cells.Item(136, 14).Value = -10644.4614
